$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 41
$ws.Range("H41").Value = 323.75
$ws.Range("I41").Value = 199
$ws.Range("J41").Value = 448.5
$ws.Range("K41").Value = 199
$ws.Range("L41").Value = 448.5
$ws.Range("M41").Value = 241
$ws.Range("N41").Value = -1328.5

# Row 53
$ws.Range("H53").Value = 248.33333
$ws.Range("I53").Value = 67
$ws.Range("J53").Value = 475
$ws.Range("K53").Value = 67
$ws.Range("L53").Value = 475
$ws.Range("M53").Value = 570

# Row 62
$ws.Range("H62").Value = 9037
$ws.Range("I62").Value = 9037
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 9037
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -8413

# Row 65
$ws.Range("H65").Value = 9037
$ws.Range("I65").Value = 9037
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 45185
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -42065

# Row 86
$ws.Range("H86").Value = 2515
$ws.Range("I86").Value = 2119.6
$ws.Range("J86").Value = 3009.25
$ws.Range("K86").Value = 2119.6
$ws.Range("L86").Value = 3009.25
$ws.Range("M86").Value = -996.5999999999999

# Row 89
$ws.Range("H89").Value = 2515
$ws.Range("I89").Value = 2119.6
$ws.Range("J89").Value = 3009.25
$ws.Range("K89").Value = 10598
$ws.Range("L89").Value = 15046.25
$ws.Range("M89").Value = -4982

# Row 107
$ws.Range("H107").Value = 1001.0476
$ws.Range("I107").Value = 1034.9445
$ws.Range("J107").Value = 797.6667
$ws.Range("K107").Value = 1034.9445
$ws.Range("L107").Value = 797.6667
$ws.Range("M107").Value = 885.0554999999999

$ws = $wb.Worksheets.Item("ARM")

# Row 110
$ws.Range("H110").Value = 2980.3845
$ws.Range("I110").Value = 1158.6364
$ws.Range("J110").Value = 13000
$ws.Range("K110").Value = 1158.6364
$ws.Range("L110").Value = 13000
$ws.Range("M110").Value = 886.3635999999999

# Row 122
$ws.Range("H122").Value = 1610
$ws.Range("I122").Value = 1218.8182
$ws.Range("J122").Value = 3044.3333
$ws.Range("K122").Value = 3656.4546
$ws.Range("L122").Value = 9132.999899999999
$ws.Range("M122").Value = -1206.4546

# Row 132
$ws.Range("H132").Value = 2982.889
$ws.Range("I132").Value = 2284.3076
$ws.Range("J132").Value = 4799.2
$ws.Range("K132").Value = 6852.9228
$ws.Range("L132").Value = 14397.6
$ws.Range("M132").Value = -4322.9228
$ws.Range("N132").Value = -19457.6

$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 2877.0454
$ws.Range("I86").Value = 1718.091
$ws.Range("J86").Value = 4036
$ws.Range("K86").Value = 1718.091
$ws.Range("L86").Value = 4036
$ws.Range("M86").Value = -595.0909999999999
$ws.Range("N86").Value = -6282

# Row 89
$ws.Range("H89").Value = 2877.0454
$ws.Range("I89").Value = 1718.091
$ws.Range("J89").Value = 4036
$ws.Range("K89").Value = 8590.455
$ws.Range("L89").Value = 20180
$ws.Range("M89").Value = -2974.455
$ws.Range("N89").Value = -31412

$ws = $wb.Worksheets.Item("CRP")

# Row 2
$ws.Range("H2").Value = 4833.3335
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 10500
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 10500
$ws.Range("M2").Value = -1887
$ws.Range("N2").Value = -10726

# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()

# Row 6
$ws.Range("H6").Value = 20000000
$ws.Range("I6").Value = 20000000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 20000000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -19999887
$ws.Range("N6").ClearContents()

# Row 10
$ws.Range("H10").Value = 1014.25
$ws.Range("I10").Value = 1014.25
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1014.25
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -875.25

# Row 12
$ws.Range("H12").Value = 27000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 27000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 27000
$ws.Range("N12").Value = -27340

# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 14
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()

# Row 16
$ws.Range("H16").Value = 2569.182
$ws.Range("I16").Value = 1904.75
$ws.Range("J16").Value = 2948.8572
$ws.Range("K16").Value = 1904.75
$ws.Range("L16").Value = 2948.8572
$ws.Range("M16").Value = -1617.75
$ws.Range("N16").Value = -3522.8572

# Row 19
$ws.Range("H19").Value = 518.3333
$ws.Range("I19").Value = 518.3333
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 518.3333
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -348.3333
$ws.Range("N19").ClearContents()

# Row 23
$ws.Range("H23").Value = 9999
$ws.Range("I23").Value = 9999
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 9999
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -9759

# Row 24
$ws.Range("H24").Value = 518.3333
$ws.Range("I24").Value = 518.3333
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 518.3333
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -348.3333
$ws.Range("N24").ClearContents()

# Row 27
$ws.Range("H27").Value = 9999
$ws.Range("I27").Value = 9999
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 9999
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -9807

# Row 113
$ws.Range("H113").Value = 2569.182
$ws.Range("I113").Value = 1904.75
$ws.Range("J113").Value = 2948.8572
$ws.Range("K113").Value = 1904.75
$ws.Range("L113").Value = 2948.8572
$ws.Range("M113").Value = 265.25
$ws.Range("N113").Value = -7288.8572

# Row 132
$ws.Range("H132").Value = 4836.727
$ws.Range("I132").Value = 4744.4
$ws.Range("J132").Value = 4913.6665
$ws.Range("K132").Value = 14233.2
$ws.Range("L132").Value = 14740.9995
$ws.Range("M132").Value = -11703.2
$ws.Range("N132").Value = -19800.9995

$ws = $wb.Worksheets.Item("CUL")

# Row 109
$ws.Range("H109").Value = 4630.8945
$ws.Range("I109").Value = 1494
$ws.Range("J109").Value = 4999.9414
$ws.Range("K109").Value = 4482
$ws.Range("L109").Value = 14999.8242
$ws.Range("M109").Value = -3442
$ws.Range("N109").Value = -17079.8242

# Row 112
$ws.Range("H112").Value = 7261.857
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 7261.857
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 21785.571
$ws.Range("N112").Value = -24001.571

$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 2583.3333
$ws.Range("I102").Value = 3375
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 3375
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -1753
$ws.Range("N102").Value = -4244

# Row 132
$ws.Range("H132").Value = 3443.1667
$ws.Range("I132").Value = 2790.5
$ws.Range("J132").Value = 4748.5
$ws.Range("K132").Value = 8371.5
$ws.Range("L132").Value = 14245.5
$ws.Range("M132").Value = -5841.5

$ws = $wb.Worksheets.Item("WVR")

# Row 132
$ws.Range("H132").Value = 1657.9166
$ws.Range("I132").Value = 1018.875
$ws.Range("J132").Value = 2936
$ws.Range("K132").Value = 3056.625
$ws.Range("L132").Value = 8808
$ws.Range("M132").Value = -526.625
$ws.Range("N132").Value = -13868
